# Fix notification enable and disable.
# Adds a "Status" column (E) with value "Done" for the existing change row,
# and narrows column D's width to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Done"
$ws.Range("E1").Value = "Status"

$ws.Columns.Item(4).ColumnWidth = 74.8

$ws.Range("E2").Select()
